# Issue #1 : Petites retouches visuelles
#
# Append two new rows to the rally-car table (Subaru / Mitsubishi, both
# from 1997), extending the worksheet's used range from A1:C5 to A1:C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Subaru - WRC Impreza GC - 1997
$ws.Range("A6").Value = "Subaru"
$ws.Range("B6").Value = "WRC Impreza GC"
# The "year" column stores its values as text elsewhere in the sheet, so
# force this cell to Text before writing the digits - otherwise Excel
# would auto-detect "1997" as a number. Re-apply the Normal style
# afterwards so the cell keeps the sheet's default formatting (no visible
# number-format change), matching the rest of the table.
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "1997"
$ws.Range("C6").Style = "Normal"

# Row 7: Mitsubishi - Lancer Evo IV - 1997
$ws.Range("A7").Value = "Mitsubishi"
$ws.Range("B7").Value = "Lancer Evo IV"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "1997"
$ws.Range("C7").Style = "Normal"
